$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.002.67"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "1.900.55"
$ws.Range("E3").Value = "  -1.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7423"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3068"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06886"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08009"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7525"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "

$ws.Range("D13").Value = "1.899.77"
$ws.Range("E13").Value = "  -1.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.240"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("E15").Value = "  -1.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.182"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.12%  "

$ws.Range("D17").Value = "30.007.23"
$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007767"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.11%  "

$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").Value = "2.152.13"
$ws.Range("E22").Value = "  -1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.084"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.335"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.98%  "

$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1266"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.040"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.355"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("E32").Value = "  -2.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.048"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05333"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.281"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7369"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01947"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.254"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4453"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.942"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.743"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8317"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.822"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "

$ws.Range("D49").Value = "2.052.53"
$ws.Range("E49").Value = "  -1.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05990"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
